$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data to insert at the top (becomes rows 2-7) and to append at the bottom (rows 28-31)
$topRows = @(
    @(0.0345138870179653, -0.0546724386513233, 0.1018617823719978),
    @(-0.0178678091615438, -0.0160352122038602, 0.0375682115554809),
    @(0.0255036242306232, -0.0097738439217209, -0.0500909499824047),
    @(0.040775254368782, 0.0166460778564214, 0.0087048299610614),
    @(0.01328631862998, 0.058643065392971, 0.0251981914043426),
    @(0.0279470849782228, 0.0377209298312664, 0.0471893399953842)
)

$tailRows = @(
    @(-0.0830776765942573, -0.1788308024406433, 0.0158824957907199),
    @(-0.0752891451120376, 0.0606283769011497, -0.0500909499824047),
    @(0.0174096599221229, -0.0296269636601209, -0.0302378293126821),
    @(0.0155770638957619, 0.0525344125926494, -0.0565050356090068)
)

$insertCount = $topRows.Count
$lastRow = $ws.UsedRange.Rows.Count

# Shift the existing data rows (2..$lastRow) down by $insertCount rows, bottom-up so
# we never overwrite a source row before it has been read. ".Value2" is used for the
# read side (reading back through ".Value" mis-serialises numerics in this host).
for ($r = $lastRow; $r -ge 2; $r--) {
    $destRow = $r + $insertCount
    for ($col = 1; $col -le 3; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $ws.Cells.Item($r, $col).Value2
    }
}

# Fill the freed-up rows (2..7) with the new top data
for ($i = 0; $i -lt $topRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $topRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $topRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $topRows[$i][2]
}

# Append the new rows after the (now shifted) existing data
$appendStart = $lastRow + $insertCount + 1
for ($i = 0; $i -lt $tailRows.Count; $i++) {
    $r = $appendStart + $i
    $ws.Cells.Item($r, 1).Value = $tailRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $tailRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $tailRows[$i][2]
}
